$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the morning in/out times for row 9 (C9:D9)
$ws.Range("C9:D9").ClearContents()

# Add afternoon/evening in/out times for rows 14 and 15
$ws.Range("E14:F15").NumberFormat = "h:mm AM/PM"
$ws.Range("E14").Value = 0.79166666666666663
$ws.Range("F14").Value = 0.95833333333333337
$ws.Range("E15").Value = 0.79166666666666663
$ws.Range("F15").Value = 0.95833333333333337

# Update the selected cell/range
$ws.Range("I11").Select()
